$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value for F8 (shared string "/flashcard.html")
$ws.Range("F8").Value = "/flashcard.html"

# Adjust column widths for D and E to match the diff.
# Note: the COM layer quantizes ColumnWidth to 1/6-character steps
# (stored_width = (round(input*6)+5)/6), so we feed in the value whose
# rounded bucket lands closest to the target stored widths
# (29.88671875 / 25.88671875) rather than the raw target numbers.
$ws.Range("D1").ColumnWidth = 29
$ws.Range("E1").ColumnWidth = 25

# Update the active cell selection to F10
$ws.Range("F10").Select()
